# Update the L5CG11 routine worksheet for Level 5, 2nd semester.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "Room" and "Block" before "Group", and
# append new "Level" and "Course" columns after it. ---
$ws.Range("H1").Value = "Room"
$ws.Range("I1").Value = "Block"
$ws.Range("J1").Value = "Group"
$ws.Range("K1").Value = "Level"
$ws.Range("L1").Value = "Course"

# --- Data rows 2-10 ---
# Columns: A Day, B Time, C Module Code, D Module Title, E Hours,
#          F Class Type, G Lecturer, H Room, I Block, J Group, K Level, L Course

$data = @(
    @("SUN","9:30-12:00","5CS022","Human Computer Interaction",2.5,"Workshop","Mr. Pravash Karki","Lab-03 Gahanapokhari","HCK","L5CG11",5,"BCS"),
    @("MON","9:30-12:00","5CS020","Distributed and Cloud Systems Programming",2.5,"Workshop","Mr. Shishir Poudel","Lab-05 Basantapur","HCK","L5CG11",5,"BCS"),
    @("TUE","9:30-11:30","5CS024","Collaborative Development",2,"Lecture","Mr. Udaya Kandel","LT-03 Walsall","WLV","L5CG(9+10+11)",5,"BCS"),
    @("TUE","12:00-14:00","5CS022","Human Computer Interaction",2,"Lecture","Mr. Pravash Karki","LT-02 Telford","WLV","L5CG(9+10+11)",5,"BCS"),
    @("WED","9:00-11:00","5CS020","Distributed and Cloud Systems Programming",2,"Lecture","Mr. Sumanta Silwal","LT-02 Telford","WLV","L5CG(9+10+11)",5,"BCS"),
    @("WED","12:00-14:00","5CS024","Collaborative Development",2,"Tutorial","Mr. Udaya Kandel","SR-02 Bilston","WLV","L5CG11",5,"BCS"),
    @("THU","9:30-12:00","5CS024","Collaborative Development",2.5,"Workshop","Mr. Udaya Kandel","SR-01 Bantok","WLV","L5CG11",5,"BCS"),
    @("FRI","9:30-11:30","5CS020","Distributed and Cloud Systems Programming",2,"Tutorial","Mr. Shishir Poudel","SR-04 Crompton","WLV","L5CG11",5,"BCS"),
    @("FRI","13:30-15:30","5CS022","Human Computer Interaction",2,"Tutorial","Mr. Pravash Karki","SR-02 Bilston","WLV","L5CG11",5,"BCS")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    $ws.Cells.Item($rowIndex, 10).Value = $row[9]
    $ws.Cells.Item($rowIndex, 11).Value = $row[10]
    $ws.Cells.Item($rowIndex, 12).Value = $row[11]
    $rowIndex++
}
